$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing BARON (col F) and Pajarito (col G) data to make
# room for the new SBB / AlphaECP columns. Old F (BARON) -> new K,
# old G (Pajarito) -> new F. Doing the copy with an explicit destination
# range preserves both the value/type and the cell style.
$ws.Range("F1").Copy($ws.Range("K1"))
$ws.Range("G1").Copy($ws.Range("F1"))
$ws.Range("F2:F5").Copy($ws.Range("K2:K5"))
$ws.Range("G2:G5").Copy($ws.Range("F2:F5"))

# --- Step 2: new header cells for SBB and AlphaECP ---
$ws.Range("G1").Value = "SBB s(gap)"
$ws.Range("H1").Value = "AlphaECP s(gap)"

$ws.Range("G1:H1").Borders.Item(9).LineStyle = 1
$ws.Range("G1:H1").Borders.Item(8).LineStyle = 1
$ws.Range("G1:H1").HorizontalAlignment = -4108

# --- Step 3: new SBB column values (text "Timed out(...)" entries) ---
$ws.Range("G2").Value = "Timed out(7%)"
$ws.Range("G3").Value = "Timed out(234%)"
$ws.Range("G4").Value = "Timed out(245%)"
$ws.Range("G5").Value = "Timed out(247%)"
$ws.Range("G5").Borders.Item(9).LineStyle = 1

# --- Step 4: new AlphaECP column values (numeric) ---
$ws.Range("H2").Value = 9
$ws.Range("H3").Value = 60
$ws.Range("H4").Value = 254
$ws.Range("H5").Value = 917

$ws.Range("H2:H4").HorizontalAlignment = -4108
$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("H5").Borders.Item(9).LineStyle = 1

# --- Step 5: footnote row ---
$ws.Range("J18").Value = "|bestbound-bestinteger|/(1e-10+|bestinteger|)"
$ws.Range("J18").Font.Name = "Courier New"
$ws.Range("J18").Font.Size = 16
$ws.Range("J18").Font.Color = 3289650
$ws.Rows.Item(18).RowHeight = 22

# --- Step 6: column widths for the re-arranged / new columns ---
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 14.333333333333332
$ws.Columns.Item(8).ColumnWidth = 13.5
$ws.Columns.Item(11).ColumnWidth = 13.333333333333332

# --- Step 7: selection / window bookkeeping ---
$ws.Range("F8").Select()

$w = $wb.Windows.Item(1)
$w.Left = 38420
$w.Top = 2280
